$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the period header text from "Aug 2022" -> "Sep 2022"
# D3 = "Aug 2022", E3 = "SD Aug 2022"
# Force text formatting first so Excel doesn't auto-convert "Sep 2022"
# into a date serial number.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "Sep 2022"
$ws.Range("E3").Value = "SD Sep 2022"

# Cell C14 (row "JUMLAH 2") currently holds a text/shared-string reference
# equal to "342,652,403"; change it to the literal number 0.
$ws.Range("C14").Value = 0
